# Applies the three paragraph-level edits described by the diff:
#   1. "3.password Encryption/decryption" -> split into 5 runs with a
#      parenthetical "(my localhost problem)" note in theme accent6 green.
#   2. " url" -> split so "url" is wrapped in a spellStart/spellEnd
#      <w:proofErr/> pair (the leading space stays its own run).
#   3. "(Bkash ,wallet)" -> split so "Bkash" gets spellStart/spellEnd and
#      the whole "Bkash ,wallet" phrase gets gramStart/gramEnd proofErr
#      markers.
#
# Word's Font/Range color setters can't express an explicit RGB value
# together with a theme color reference at the same time, and there is
# no supported way to splice a bare <w:proofErr/> marker between two
# runs of an *existing* paragraph through the Range/Selection API - so
# each affected paragraph is rebuilt wholesale via Range.InsertXML with
# literal OOXML that reproduces the target run/proofErr structure
# exactly, which lets both concerns be expressed precisely in one shot.

$d = $word.ActiveDocument

function Replace-ParagraphXml($searchText, $bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText)
    if (-not $found) {
        throw "Could not find paragraph text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $prange = $para.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $prange.InsertXML($xml)
}

# --- Change 1: "3.password Encryption/decryption" paragraph -------------
$body1 = '<w:p w14:paraId="1BC2EDBE" w14:textId="77777777" w:rsidR="00437065" w:rsidRPr="003C21D0" w:rsidRDefault="00437065" w:rsidP="00437065">' +
    '<w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="003C21D0"><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>3.password Encryption/</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>decryption</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/><w:lang w:val="en-GB"/></w:rPr><w:t>my localhost problem</w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
Replace-ParagraphXml "3.password Encryption/decryption" $body1

# --- Change 2: " url" -> proofErr-wrapped "url" --------------------------
$body2 = '<w:p w14:paraId="7B381F55" w14:textId="77777777" w:rsidR="00437065" w:rsidRDefault="00437065" w:rsidP="00437065">' +
    '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">5.front </w:t></w:r>' +
    '<w:r w:rsidR="003C21D0"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>SEO</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> friendly</w:t></w:r>' +
    '<w:r w:rsidR="003C21D0"><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>url</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Replace-ParagraphXml "5.front SEO friendly url" $body2

# --- Change 3: "(Bkash ,wallet)" -> proofErr-wrapped phrase --------------
$body3 = '<w:p w14:paraId="39BF7B07" w14:textId="77777777" w:rsidR="00D57493" w:rsidRPr="00D57493" w:rsidRDefault="00D57493" w:rsidP="009A28A3">' +
    '<w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00D57493"><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>Bkash</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> ,wallet</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
Replace-ParagraphXml "(Bkash ,wallet)" $body3

Write-Host "Applied all 3 paragraph edits."
